# Regenerate the localization-status report after handback:
# the two source files (62e2ac14... and 5803e10a...) swap report rows and
# move from "Ready for handoff" to "Handed back: in sync with en-US", with
# fresh handback timestamps for both the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value2 = "5803e10a-454f-4e1e-9572-eef29b5960e9.md"
$ws1.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws1.Range("C2").Value2 = "Handed back: in sync with en-US"

$ws1.Range("A3").Value2 = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"
$ws1.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws1.Range("C3").Value2 = "Handed back: in sync with en-US"

$h1a = $ws1.Hyperlinks.Item(1)
$h1a.Address = "https://github.com/OpenLocalizationTest/oltest/blob/9fe253d9c8681c2154dc7707dc3a9b1e14cb448c/e2e/5803e10a-454f-4e1e-9572-eef29b5960e9.md"
$h1a.TextToDisplay = "5803e10a-454f-4e1e-9572-eef29b5960e9.md"

$h1b = $ws1.Hyperlinks.Item(2)
$h1b.Address = "https://github.com/OpenLocalizationTest/oltest/blob/9fe253d9c8681c2154dc7707dc3a9b1e14cb448c/e2e/62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"
$h1b.TextToDisplay = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"

# ---- zh-cn sheet ----
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value2 = "5803e10a-454f-4e1e-9572-eef29b5960e9.md"
$ws2.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws2.Range("C2").Value2 = "5803e10a-454f-4e1e-9572-eef29b5960e9.208bfd7be456e4feaa124079e4b95d71bf8108a1.zh-cn.xlf"
$ws2.Range("D2").Value2 = "2016-03-08 07:07:16"
$ws2.Range("E2").Value2 = "5803e10a-454f-4e1e-9572-eef29b5960e9.md"
$ws2.Range("F2").Value2 = "5803e10a-454f-4e1e-9572-eef29b5960e9.208bfd7be456e4feaa124079e4b95d71bf8108a1.zh-cn.xlf"
$ws2.Range("G2").Value2 = "2016-03-08 07:08:00"
$ws2.Range("H2").Value2 = "Include"

$ws2.Range("A3").Value2 = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"
$ws2.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws2.Range("C3").Value2 = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.3472016d53bb0fd6d29145dff96b6f8971dfae1b.zh-cn.xlf"
$ws2.Range("D3").Value2 = "2016-03-08 07:07:16"
$ws2.Range("E3").Value2 = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"
$ws2.Range("F3").Value2 = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.3472016d53bb0fd6d29145dff96b6f8971dfae1b.zh-cn.xlf"
$ws2.Range("G3").Value2 = "2016-03-08 07:08:00"
$ws2.Range("H3").Value2 = "Include"

$h2a = $ws2.Hyperlinks.Item(1)
$h2a.Address = "https://github.com/OpenLocalizationTest/oltest/blob/9fe253d9c8681c2154dc7707dc3a9b1e14cb448c/e2e/5803e10a-454f-4e1e-9572-eef29b5960e9.md"
$h2a.TextToDisplay = "5803e10a-454f-4e1e-9572-eef29b5960e9.md"

$h2b = $ws2.Hyperlinks.Item(2)
$h2b.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46cc4b6cd022077f4518774a198c3cfeec99149b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/mt/5803e10a-454f-4e1e-9572-eef29b5960e9.208bfd7be456e4feaa124079e4b95d71bf8108a1.zh-cn.xlf"
$h2b.TextToDisplay = "5803e10a-454f-4e1e-9572-eef29b5960e9.208bfd7be456e4feaa124079e4b95d71bf8108a1.zh-cn.xlf"

$h2c = $ws2.Hyperlinks.Item(3)
$h2c.Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c5988671ab236c48d1308b0a7f30f2420d52bb86/e2e/5803e10a-454f-4e1e-9572-eef29b5960e9.md"
$h2c.TextToDisplay = "5803e10a-454f-4e1e-9572-eef29b5960e9.md"

$h2d = $ws2.Hyperlinks.Item(4)
$h2d.Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5ab91f02d49feda957d1e12a8b500d50e03b65d6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5803e10a-454f-4e1e-9572-eef29b5960e9.208bfd7be456e4feaa124079e4b95d71bf8108a1.zh-cn.xlf"
$h2d.TextToDisplay = "5803e10a-454f-4e1e-9572-eef29b5960e9.208bfd7be456e4feaa124079e4b95d71bf8108a1.zh-cn.xlf"

$h2e = $ws2.Hyperlinks.Item(5)
$h2e.Address = "https://github.com/OpenLocalizationTest/oltest/blob/9fe253d9c8681c2154dc7707dc3a9b1e14cb448c/e2e/62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"
$h2e.TextToDisplay = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"

$h2f = $ws2.Hyperlinks.Item(6)
$h2f.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46cc4b6cd022077f4518774a198c3cfeec99149b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/mt/62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.3472016d53bb0fd6d29145dff96b6f8971dfae1b.zh-cn.xlf"
$h2f.TextToDisplay = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.3472016d53bb0fd6d29145dff96b6f8971dfae1b.zh-cn.xlf"

$h2g = $ws2.Hyperlinks.Item(7)
$h2g.Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c5988671ab236c48d1308b0a7f30f2420d52bb86/e2e/62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"
$h2g.TextToDisplay = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"

$h2h = $ws2.Hyperlinks.Item(8)
$h2h.Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5ab91f02d49feda957d1e12a8b500d50e03b65d6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.3472016d53bb0fd6d29145dff96b6f8971dfae1b.zh-cn.xlf"
$h2h.TextToDisplay = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.3472016d53bb0fd6d29145dff96b6f8971dfae1b.zh-cn.xlf"

# ---- de-de sheet ----
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value2 = "5803e10a-454f-4e1e-9572-eef29b5960e9.md"
$ws3.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws3.Range("C2").Value2 = "5803e10a-454f-4e1e-9572-eef29b5960e9.208bfd7be456e4feaa124079e4b95d71bf8108a1.de-de.xlf"
$ws3.Range("D2").Value2 = "2016-03-08 07:07:27"
$ws3.Range("E2").Value2 = "5803e10a-454f-4e1e-9572-eef29b5960e9.md"
$ws3.Range("F2").Value2 = "5803e10a-454f-4e1e-9572-eef29b5960e9.208bfd7be456e4feaa124079e4b95d71bf8108a1.de-de.xlf"
$ws3.Range("G2").Value2 = "2016-03-08 07:08:17"
$ws3.Range("H2").Value2 = "Include"

$ws3.Range("A3").Value2 = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"
$ws3.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws3.Range("C3").Value2 = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.3472016d53bb0fd6d29145dff96b6f8971dfae1b.de-de.xlf"
$ws3.Range("D3").Value2 = "2016-03-08 07:07:27"
$ws3.Range("E3").Value2 = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"
$ws3.Range("F3").Value2 = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.3472016d53bb0fd6d29145dff96b6f8971dfae1b.de-de.xlf"
$ws3.Range("G3").Value2 = "2016-03-08 07:08:17"
$ws3.Range("H3").Value2 = "Include"

$h3a = $ws3.Hyperlinks.Item(1)
$h3a.Address = "https://github.com/OpenLocalizationTest/oltest/blob/9fe253d9c8681c2154dc7707dc3a9b1e14cb448c/e2e/5803e10a-454f-4e1e-9572-eef29b5960e9.md"
$h3a.TextToDisplay = "5803e10a-454f-4e1e-9572-eef29b5960e9.md"

$h3b = $ws3.Hyperlinks.Item(2)
$h3b.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fddc81e046fe6cb3e85733f9f9bb27c8c9aef1ee/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/mt/5803e10a-454f-4e1e-9572-eef29b5960e9.208bfd7be456e4feaa124079e4b95d71bf8108a1.de-de.xlf"
$h3b.TextToDisplay = "5803e10a-454f-4e1e-9572-eef29b5960e9.208bfd7be456e4feaa124079e4b95d71bf8108a1.de-de.xlf"

$h3c = $ws3.Hyperlinks.Item(3)
$h3c.Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/16155aa1e8be1dd6ec20370d6658b051b4c9e7d2/e2e/5803e10a-454f-4e1e-9572-eef29b5960e9.md"
$h3c.TextToDisplay = "5803e10a-454f-4e1e-9572-eef29b5960e9.md"

$h3d = $ws3.Hyperlinks.Item(4)
$h3d.Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a859616c4ef4872c2d98c6600a80893ed799faf6/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5803e10a-454f-4e1e-9572-eef29b5960e9.208bfd7be456e4feaa124079e4b95d71bf8108a1.de-de.xlf"
$h3d.TextToDisplay = "5803e10a-454f-4e1e-9572-eef29b5960e9.208bfd7be456e4feaa124079e4b95d71bf8108a1.de-de.xlf"

$h3e = $ws3.Hyperlinks.Item(5)
$h3e.Address = "https://github.com/OpenLocalizationTest/oltest/blob/9fe253d9c8681c2154dc7707dc3a9b1e14cb448c/e2e/62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"
$h3e.TextToDisplay = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"

$h3f = $ws3.Hyperlinks.Item(6)
$h3f.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fddc81e046fe6cb3e85733f9f9bb27c8c9aef1ee/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/mt/62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.3472016d53bb0fd6d29145dff96b6f8971dfae1b.de-de.xlf"
$h3f.TextToDisplay = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.3472016d53bb0fd6d29145dff96b6f8971dfae1b.de-de.xlf"

$h3g = $ws3.Hyperlinks.Item(7)
$h3g.Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/16155aa1e8be1dd6ec20370d6658b051b4c9e7d2/e2e/62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"
$h3g.TextToDisplay = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.md"

$h3h = $ws3.Hyperlinks.Item(8)
$h3h.Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a859616c4ef4872c2d98c6600a80893ed799faf6/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.3472016d53bb0fd6d29145dff96b6f8971dfae1b.de-de.xlf"
$h3h.TextToDisplay = "62e2ac14-30fe-4f47-ab8c-0ed8657f6e84.3472016d53bb0fd6d29145dff96b6f8971dfae1b.de-de.xlf"

Write-Output "Localization status report regenerated for handback."
